$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "23.798.14"
Set-TextValue "E2" "  -0.75%  "
Set-TextValue "D3" "1.636.91"
Set-TextValue "E3" "  -1.09%  "
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "308.87"
Set-TextValue "E5" "  -0.22%  "
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.09%  "
Set-TextValue "D7" "0.3867"
Set-TextValue "E7" "  -1.20%  "
Set-TextValue "D8" "0.3802"
Set-TextValue "E8" "  -2.14%  "
Set-TextValue "D9" "50.48"
Set-TextValue "E9" "  -2.51%  "
Set-TextValue "D10" "1.322"
Set-TextValue "E10" "  -3.57%  "
Set-TextValue "D11" "1.002"
Set-TextValue "E11" "  -0.13%  "
Set-TextValue "D12" "0.08363"
Set-TextValue "E12" "  -1.57%  "
Set-TextValue "D13" "23.58"
Set-TextValue "E13" "  -2.42%  "
Set-TextValue "D14" "6.949"
Set-TextValue "E14" "  -4.37%  "
Set-TextValue "D15" "7.810"
Set-TextValue "E15" "  -3.68%  "
Set-TextValue "D16" "0.00001302"
Set-TextValue "E16" "  -1.20%  "
Set-TextValue "D17" "1.640.39"
Set-TextValue "E17" "  -0.88%  "
Set-TextValue "D18" "93.33"
Set-TextValue "E18" "  -1.97%  "
Set-TextValue "D19" "0.06944"
Set-TextValue "E19" "  -0.33%  "
Set-TextValue "D20" "19.34"
Set-TextValue "E20" "  -3.29%  "
Set-TextValue "D21" "6.838"
Set-TextValue "E21" "  -2.18%  "
Set-TextValue "E22" "  +0.07%  "
Set-TextValue "D23" "13.51"
Set-TextValue "E23" "  -1.59%  "
Set-TextValue "D24" "23.811.74"
Set-TextValue "E24" "  -0.73%  "
Set-TextValue "D25" "2.435"
Set-TextValue "E25" "  -2.63%  "
Set-TextValue "D26" "2.860"
Set-TextValue "E26" "  -9.29%  "
Set-TextValue "D27" "21.73"
Set-TextValue "E27" "  -2.58%  "
Set-TextValue "D28" "153.48"
Set-TextValue "E28" "  -0.27%  "
Set-TextValue "D29" "5.503"
Set-TextValue "E29" "  +4.01%  "
Set-TextValue "D30" "136.08"
Set-TextValue "E30" "  -2.87%  "
Set-TextValue "D31" "7.805"
Set-TextValue "E31" "  -0.87%  "
Set-TextValue "D32" "2.495"
Set-TextValue "E32" "  +0.60%  "
Set-TextValue "D33" "1.821.21"
Set-TextValue "D34" "0.07936"
Set-TextValue "E34" "  -2.33%  "
Set-TextValue "D35" "0.9754"
Set-TextValue "E35" "  -6.75%  "
Set-TextValue "D36" "0.02888"
Set-TextValue "E36" "  -4.24%  "
Set-TextValue "D37" "6.567"
Set-TextValue "E37" "  -1.94%  "
Set-TextValue "D38" "0.2646"
Set-TextValue "E38" "  -2.30%  "
Set-TextValue "D39" "10.36"
Set-TextValue "E39" "  -7.49%  "
Set-TextValue "D40" "0.09055"
Set-TextValue "E40" "  -1.09%  "
Set-TextValue "D41" "0.7449"
Set-TextValue "E41" "  -2.04%  "
Set-TextValue "D42" "13.19"
Set-TextValue "E42" "  -2.46%  "
Set-TextValue "E43" "  -0.83%  "
Set-TextValue "D44" "16.49"
Set-TextValue "E44" "  -1.18%  "
Set-TextValue "D45" "0.6853"
Set-TextValue "E45" "  -2.64%  "
Set-TextValue "D46" "2.402"
Set-TextValue "E46" "  -4.10%  "
Set-TextValue "D47" "4.067"
Set-TextValue "E47" "  -0.46%  "
Set-TextValue "D48" "1.001"
Set-TextValue "E48" "  +0.07%  "
Set-TextValue "E49" "  -1.89%  "
Set-TextValue "D50" "133.86"
Set-TextValue "E50" "  -0.99%  "
Set-TextValue "D51" "1.211"

Write-Output "Done applying 94 cell updates"
